$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Athaya Devin Argyadama ---
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Range("A6").Value = 44386.350098888885
$ws.Range("B6").Value = "dvndevin05@gmail.com"
$ws.Range("C6").Value = "Athaya Devin Argyadama"
$ws.Range("D6").Value = "X-IPA 1"
$ws.Range("E6").Value = "https://drive.google.com/open?id=1Y3C1IELh8uIDZDZGV5vyLed_quhIV1bI"
$ws.Hyperlinks.Add($ws.Range("E6"), "https://drive.google.com/open?id=1Y3C1IELh8uIDZDZGV5vyLed_quhIV1bI")
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# --- Row 7: Alfyan Cana Dwi Fakhrudi ---
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("A7").Value = 44387.62598138889
$ws.Range("B7").Value = "alfyanfakhrudi@gmail.com"
$ws.Range("C7").Value = "Alfyan Cana Dwi Fakhrudi"
$ws.Range("D7").Value = "X IPA 4"
$ws.Range("E7").Value = "https://drive.google.com/open?id=1l77MP4jDSCYw1bfzMKGJleJ99AoSJ1VC"
$ws.Hyperlinks.Add($ws.Range("E7"), "https://drive.google.com/open?id=1l77MP4jDSCYw1bfzMKGJleJ99AoSJ1VC")
$wb.Styles.Item("Hyperlink").Delete()
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)
